$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly resample of the daily "Fruta, Terminal La Palmera de La Serena - Coco"
# data: each data row (2..30) is reassigned the Fecha/Volumen/Precio fields
# that used to belong to a different row in the original daily table, as
# part of switching the consolidated subconjunto to a weekly ("semanal")
# sampling logic.

$rowData = @(
    [PSCustomObject]@{ Row=2; Year=2021; Month=5; Day=20; M=100; N=19500; O=20000; P=19750; S=988 }
    [PSCustomObject]@{ Row=3; Year=2021; Month=8; Day=23; M=160; N=21000; O=22000; P=21500; S=1075 }
    [PSCustomObject]@{ Row=4; Year=2021; Month=6; Day=17; M=140; N=20000; O=21000; P=20500; S=1025 }
    [PSCustomObject]@{ Row=5; Year=2021; Month=8; Day=10; M=200; N=20000; O=21000; P=20500; S=1025 }
    [PSCustomObject]@{ Row=6; Year=2021; Month=10; Day=5; M=200; N=19000; O=20000; P=19500; S=975 }
    [PSCustomObject]@{ Row=7; Year=2021; Month=4; Day=29; M=100; N=20000; O=21000; P=20500; S=1025 }
    [PSCustomObject]@{ Row=8; Year=2021; Month=9; Day=9; M=100; N=20000; O=21000; P=20500; S=1025 }
    [PSCustomObject]@{ Row=9; Year=2021; Month=8; Day=9; M=160; N=20000; O=21000; P=20500; S=1025 }
    [PSCustomObject]@{ Row=10; Year=2021; Month=6; Day=3; M=160; N=19000; O=20000; P=19500; S=975 }
    [PSCustomObject]@{ Row=11; Year=2021; Month=8; Day=12; M=160; N=20000; O=21000; P=20500; S=1025 }
    [PSCustomObject]@{ Row=12; Year=2021; Month=7; Day=30; M=160; N=20000; O=21000; P=20500; S=1025 }
    [PSCustomObject]@{ Row=13; Year=2021; Month=9; Day=6; M=160; N=20000; O=21000; P=20500; S=1025 }
    [PSCustomObject]@{ Row=14; Year=2021; Month=5; Day=17; M=100; N=19500; O=20000; P=19750; S=988 }
    [PSCustomObject]@{ Row=15; Year=2021; Month=9; Day=3; M=140; N=20000; O=21000; P=20500; S=1025 }
    [PSCustomObject]@{ Row=16; Year=2021; Month=9; Day=23; M=100; N=19500; O=20000; P=19750; S=988 }
    [PSCustomObject]@{ Row=17; Year=2021; Month=6; Day=18; M=100; N=20000; O=21000; P=20500; S=1025 }
    [PSCustomObject]@{ Row=18; Year=2021; Month=8; Day=26; M=100; N=20000; O=21000; P=20500; S=1025 }
    [PSCustomObject]@{ Row=19; Year=2021; Month=4; Day=15; M=100; N=18000; O=19000; P=18500; S=925 }
    [PSCustomObject]@{ Row=20; Year=2021; Month=5; Day=10; M=160; N=19500; O=20000; P=19750; S=988 }
    [PSCustomObject]@{ Row=21; Year=2021; Month=10; Day=4; M=40; N=19500; O=20000; P=19750; S=988 }
    [PSCustomObject]@{ Row=22; Year=2021; Month=8; Day=20; M=100; N=20000; O=21000; P=20500; S=1025 }
    [PSCustomObject]@{ Row=23; Year=2021; Month=8; Day=19; M=200; N=20000; O=21000; P=20500; S=1025 }
    [PSCustomObject]@{ Row=24; Year=2021; Month=9; Day=2; M=160; N=20000; O=21000; P=20500; S=1025 }
    [PSCustomObject]@{ Row=25; Year=2021; Month=5; Day=19; M=200; N=19000; O=20000; P=19500; S=975 }
    [PSCustomObject]@{ Row=26; Year=2021; Month=8; Day=2; M=200; N=20000; O=21000; P=20500; S=1025 }
    [PSCustomObject]@{ Row=27; Year=2021; Month=5; Day=27; M=100; N=19500; O=20000; P=19750; S=988 }
    [PSCustomObject]@{ Row=28; Year=2021; Month=8; Day=27; M=260; N=20000; O=22000; P=21115; S=1056 }
    [PSCustomObject]@{ Row=29; Year=2021; Month=9; Day=28; M=200; N=20000; O=21000; P=20500; S=1025 }
    [PSCustomObject]@{ Row=30; Year=2021; Month=9; Day=27; M=100; N=20000; O=21000; P=20500; S=1025 }
)

foreach ($r in $rowData) {
    $ws.Cells.Item($r.Row, 4).Value = Get-Date -Year $r.Year -Month $r.Month -Day $r.Day -Hour 0 -Minute 0 -Second 0
    $ws.Cells.Item($r.Row, 13).Value = $r.M
    $ws.Cells.Item($r.Row, 14).Value = $r.N
    $ws.Cells.Item($r.Row, 15).Value = $r.O
    $ws.Cells.Item($r.Row, 16).Value = $r.P
    $ws.Cells.Item($r.Row, 19).Value = $r.S
}
